# Auto-generated edit script: updates cryptos price/volume data (cap3k GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''24.642.18'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '''1.675.46'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.74%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''313.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '''0.3933'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("D8").Value = '''0.3952'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.91%  '
$ws.Range("D9").Value = '''1.002'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("D10").Value = '''1.404'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.56%  '
$ws.Range("D11").Value = '''50.99'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -5.18%  '
$ws.Range("D12").Value = '''0.08653'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").Value = '''25.22'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.42%  '
$ws.Range("D14").Value = '''7.340'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = '''0.00001316'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.47%  '
$ws.Range("D16").Value = '''7.717'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -4.01%  '
$ws.Range("D17").Value = '''1.681.46'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.12%  '
$ws.Range("D18").Value = '''94.05'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").Value = '''0.07017'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.34%  '
$ws.Range("D20").Value = '''21.10'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").Value = '''7.083'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").Value = '''1.001'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '''13.94'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.03%  '
$ws.Range("D24").Value = '''24.649.33'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").Value = '''2.351'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").Value = '''2.780'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.40%  '
$ws.Range("D27").Value = '''23.10'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").Value = '''5.881'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -8.46%  '
$ws.Range("D29").Value = '''159.88'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.25%  '
$ws.Range("D30").Value = '''146.49'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("D31").Value = '''8.389'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("D32").Value = '''2.496'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +10.01%  '
$ws.Range("D33").Value = '''1.864.60'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("D34").Value = '''0.03088'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.16%  '
$ws.Range("D35").Value = '''0.08303'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.09%  '
$ws.Range("D36").Value = '''6.973'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -5.42%  '
$ws.Range("D37").Value = '''0.2803'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.22%  '
$ws.Range("D38").Value = '''0.9903'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.14%  '
$ws.Range("D39").Value = '''0.09642'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("D40").Value = '''1.527'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.54%  '
$ws.Range("D41").Value = '''10.30'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.14%  '
$ws.Range("D42").Value = '''0.7916'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.18%  '
$ws.Range("D43").Value = '''13.52'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.03%  '
$ws.Range("D44").Value = '''16.66'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.86%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''2.563'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.72%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.7106'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.97%  '
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D49").Value = '''1.001'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.03%  '
$ws.Range("D50").Value = '''1.327'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.47%  '
$ws.Range("D51").Value = '''137.95'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.31%  '
